$wb = $excel.ActiveWorkbook
$inv = $wb.Worksheets.Item("Inventory")

# --- Add the new "Boards" worksheet right after "Inventory" ---
$ws = $wb.Worksheets.Add($null, $inv)
$ws.Name = "Boards"

# --- Column widths (best effort; engine quantizes to 1/7 char) ---
$ws.Columns.Item(1).ColumnWidth = 29.375558035714285
$ws.Columns.Item(2).ColumnWidth = 9.285714285714286
$ws.Columns.Item(3).ColumnWidth = 9.857142857142858
$ws.Columns.Item(4).ColumnWidth = 9.857142857142858
$ws.Columns.Item(5).ColumnWidth = 28.714285714285715

# --- Row heights that differ from the default ---
$ws.Rows.Item(1).RowHeight = 13
$ws.Rows.Item(5).RowHeight = 13

# ============================================================
# Cell values are written in a very specific order below so the
# shared-string table is built up in the same sequence as the
# target workbook (new entries are appended to sst in
# first-seen order).
# ============================================================

# Title
$ws.Range("A1").Value = "Inventory of evaluation boards"

# Header row (write C5/D5 before A5 so "Available"/"Assembled" get
# lower shared-string indices than "Name")
$ws.Range("C5").Value = "Available"
$ws.Range("D5").Value = "Assembled"
$ws.Range("A5").Value = "Name"
$ws.Range("B5").Value = "Cnt"
$ws.Range("E5").Value = "Notes"

# Names (pre-sort insertion order -- alphabetical list entered
# first, then "microstrip..." appended last, matching how the
# shared strings were interned before the sheet got sorted)
$ws.Range("A6").Value = "amp-LNA-L-band-TQP3M9037"
$ws.Range("A7").Value = "amp-LNA-S-band-TQP3M9037"
$ws.Range("A8").Value = "amp-LNA-UHF-TQP3M9036"
$ws.Range("A9").Value = "bpf-IF-SAW-856930"
$ws.Range("A11").Value = "mixer-MAC-24+"
$ws.Range("A12").Value = "power-linear-LP38692"
$ws.Range("A13").Value = "power-switching-TPS6211x"
$ws.Range("A14").Value = "synthesizer-RF-Si4123"
$ws.Range("A10").Value = "microstrip-test-4-layer-sma-v1"

# Status values, first occurrence order: RCVD, PEND, ORD ?
$ws.Range("B14").Value = "RCVD"
$ws.Range("B9").Value = "PEND"
$ws.Range("B6").Value = "ORD ?"

# Remaining duplicate status cells
$ws.Range("B10").Value = "RCVD"
$ws.Range("B11").Value = "PEND"
$ws.Range("B12").Value = "PEND"
$ws.Range("B13").Value = "PEND"
$ws.Range("B7").Value = "ORD ?"
$ws.Range("B8").Value = "ORD ?"

# Numeric counts
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1

# ============================================================
# Styling
# ============================================================

# A1 title -> bold (same style as Inventory!A1)
$ws.Range("A1").Font.Bold = $true

# B1 -> bold + centered (new style)
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108

# A5 "Name" header -> bold + bottom border (same as Inventory!C4)
$inv.Range("C4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# B5:E5 headers -> bold + bottom border + centered (same as Inventory!A4)
$inv.Range("A4").Copy()
$ws.Range("B5:E5").PasteSpecial(-4122)

# Status column cells -> centered (same style as Inventory!A5)
$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("B9").HorizontalAlignment = -4108
$ws.Range("B11").HorizontalAlignment = -4108
$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("B13").HorizontalAlignment = -4108

# RCVD cells -> green fill + centered (new style)
$ws.Range("B10").Interior.Color = 5296274
$ws.Range("B10").HorizontalAlignment = -4108
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# Numeric count cells -> centered
$ws.Range("C10").HorizontalAlignment = -4108
$ws.Range("D10").HorizontalAlignment = -4108
$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("D14").HorizontalAlignment = -4108

# --- Sheet view: zoom 120%, frozen none, selection A20/A19:A20 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 120
$ws.Range("A20").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Restore Inventory as the active sheet/selection ---
$inv.Activate()
$inv.Range("C38").Select()
